$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = '{''label'': ''LABEL_0'', ''score'': 0.7060643434524536}'
$ws.Range("H3").Value = '{''label'': ''LABEL_0'', ''score'': 0.7031141519546509}'
$ws.Range("H4").Value = '{''label'': ''LABEL_0'', ''score'': 0.7124890685081482}'
$ws.Range("H5").Value = '{''label'': ''LABEL_0'', ''score'': 0.6995008587837219}'
$ws.Range("H6").Value = '{''label'': ''LABEL_0'', ''score'': 0.7029309272766113}'
$ws.Range("H7").Value = '{''label'': ''LABEL_0'', ''score'': 0.7012220025062561}'
$ws.Range("H8").Value = '{''label'': ''LABEL_0'', ''score'': 0.7027474641799927}'
$ws.Range("H9").Value = '{''label'': ''LABEL_0'', ''score'': 0.6979803442955017}'
$ws.Range("H10").Value = '{''label'': ''LABEL_0'', ''score'': 0.7009510397911072}'
$ws.Range("H11").Value = '{''label'': ''LABEL_0'', ''score'': 0.7117259502410889}'
$ws.Range("H12").Value = '{''label'': ''LABEL_0'', ''score'': 0.7162072658538818}'
$ws.Range("H13").Value = '{''label'': ''LABEL_0'', ''score'': 0.7009451985359192}'
$ws.Range("H14").Value = '{''label'': ''LABEL_0'', ''score'': 0.7081892490386963}'
$ws.Range("H15").Value = '{''label'': ''LABEL_0'', ''score'': 0.7122577428817749}'
$ws.Range("H16").Value = '{''label'': ''LABEL_0'', ''score'': 0.694705069065094}'
$ws.Range("H17").Value = '{''label'': ''LABEL_0'', ''score'': 0.6902344822883606}'
$ws.Range("H18").Value = '{''label'': ''LABEL_0'', ''score'': 0.7087001204490662}'
$ws.Range("H19").Value = '{''label'': ''LABEL_0'', ''score'': 0.7044041156768799}'
$ws.Range("H20").Value = '{''label'': ''LABEL_0'', ''score'': 0.7215670347213745}'
$ws.Range("H21").Value = '{''label'': ''LABEL_0'', ''score'': 0.7173374891281128}'
$ws.Range("H22").Value = '{''label'': ''LABEL_0'', ''score'': 0.6990979909896851}'
$ws.Range("H23").Value = '{''label'': ''LABEL_0'', ''score'': 0.7108794450759888}'
$ws.Range("H24").Value = '{''label'': ''LABEL_0'', ''score'': 0.6965128183364868}'
$ws.Range("H25").Value = '{''label'': ''LABEL_0'', ''score'': 0.7174103856086731}'
$ws.Range("H26").Value = '{''label'': ''LABEL_0'', ''score'': 0.714772641658783}'
$ws.Range("H27").Value = '{''label'': ''LABEL_0'', ''score'': 0.714772641658783}'
$ws.Range("H28").Value = '{''label'': ''LABEL_0'', ''score'': 0.6991731524467468}'
$ws.Range("H29").Value = '{''label'': ''LABEL_0'', ''score'': 0.7042048573493958}'
$ws.Range("H30").Value = '{''label'': ''LABEL_0'', ''score'': 0.7437578439712524}'
$ws.Range("H31").Value = '{''label'': ''LABEL_0'', ''score'': 0.7018510103225708}'
$ws.Range("H32").Value = '{''label'': ''LABEL_0'', ''score'': 0.712437629699707}'
$ws.Range("H33").Value = '{''label'': ''LABEL_0'', ''score'': 0.712437629699707}'
$ws.Range("H34").Value = '{''label'': ''LABEL_0'', ''score'': 0.712437629699707}'
$ws.Range("H35").Value = '{''label'': ''LABEL_0'', ''score'': 0.7006770372390747}'
$ws.Range("H36").Value = '{''label'': ''LABEL_0'', ''score'': 0.6990013718605042}'
$ws.Range("H37").Value = '{''label'': ''LABEL_0'', ''score'': 0.7057111263275146}'
$ws.Range("H38").Value = '{''label'': ''LABEL_0'', ''score'': 0.6992059946060181}'
$ws.Range("H39").Value = '{''label'': ''LABEL_0'', ''score'': 0.6891406178474426}'
$ws.Range("H40").Value = '{''label'': ''LABEL_0'', ''score'': 0.6955220699310303}'
$ws.Range("H41").Value = '{''label'': ''LABEL_0'', ''score'': 0.7196774482727051}'
$ws.Range("H42").Value = '{''label'': ''LABEL_0'', ''score'': 0.7043949365615845}'
$ws.Range("H43").Value = '{''label'': ''LABEL_0'', ''score'': 0.6948150396347046}'
$ws.Range("H44").Value = '{''label'': ''LABEL_0'', ''score'': 0.6972740292549133}'
$ws.Range("H45").Value = '{''label'': ''LABEL_0'', ''score'': 0.6998518705368042}'
$ws.Range("H46").Value = '{''label'': ''LABEL_0'', ''score'': 0.6952220797538757}'
$ws.Range("H47").Value = '{''label'': ''LABEL_0'', ''score'': 0.7081985473632812}'
$ws.Range("H48").Value = '{''label'': ''LABEL_0'', ''score'': 0.7076496481895447}'
$ws.Range("H49").Value = '{''label'': ''LABEL_0'', ''score'': 0.7043262720108032}'
$ws.Range("H50").Value = '{''label'': ''LABEL_0'', ''score'': 0.7235633134841919}'
$ws.Range("H51").Value = '{''label'': ''LABEL_0'', ''score'': 0.7047570943832397}'
$ws.Range("H52").Value = '{''label'': ''LABEL_0'', ''score'': 0.693343997001648}'
$ws.Range("H53").Value = '{''label'': ''LABEL_0'', ''score'': 0.711338222026825}'
$ws.Range("H54").Value = '{''label'': ''LABEL_0'', ''score'': 0.7013468742370605}'
$ws.Range("H55").Value = '{''label'': ''LABEL_0'', ''score'': 0.6828898787498474}'
$ws.Range("H56").Value = '{''label'': ''LABEL_0'', ''score'': 0.6985781788825989}'
$ws.Range("H57").Value = '{''label'': ''LABEL_0'', ''score'': 0.7072131037712097}'
$ws.Range("H58").Value = '{''label'': ''LABEL_0'', ''score'': 0.7100904583930969}'
$ws.Range("H59").Value = '{''label'': ''LABEL_0'', ''score'': 0.7125943303108215}'
$ws.Range("H60").Value = '{''label'': ''LABEL_0'', ''score'': 0.7129694223403931}'
$ws.Range("H61").Value = '{''label'': ''LABEL_0'', ''score'': 0.6952589154243469}'
$ws.Range("H62").Value = '{''label'': ''LABEL_0'', ''score'': 0.7100904583930969}'
$ws.Range("H63").Value = '{''label'': ''LABEL_0'', ''score'': 0.7134063243865967}'
$ws.Range("H64").Value = '{''label'': ''LABEL_0'', ''score'': 0.7165772318840027}'
$ws.Range("H65").Value = '{''label'': ''LABEL_0'', ''score'': 0.7070124745368958}'
$ws.Range("H66").Value = '{''label'': ''LABEL_0'', ''score'': 0.7149694561958313}'
$ws.Range("H67").Value = '{''label'': ''LABEL_0'', ''score'': 0.7006003260612488}'
$ws.Range("H68").Value = '{''label'': ''LABEL_0'', ''score'': 0.6897069811820984}'
$ws.Range("H69").Value = '{''label'': ''LABEL_0'', ''score'': 0.7066183090209961}'
$ws.Range("H70").Value = '{''label'': ''LABEL_0'', ''score'': 0.6868455410003662}'
$ws.Range("H71").Value = '{''label'': ''LABEL_0'', ''score'': 0.7175314426422119}'
$ws.Range("H72").Value = '{''label'': ''LABEL_0'', ''score'': 0.7129970192909241}'
$ws.Range("H73").Value = '{''label'': ''LABEL_0'', ''score'': 0.7049190998077393}'
$ws.Range("H74").Value = '{''label'': ''LABEL_0'', ''score'': 0.6992518305778503}'
$ws.Range("H75").Value = '{''label'': ''LABEL_0'', ''score'': 0.6826577186584473}'
$ws.Range("H76").Value = '{''label'': ''LABEL_0'', ''score'': 0.7148903012275696}'
$ws.Range("H77").Value = '{''label'': ''LABEL_0'', ''score'': 0.7103974223136902}'
$ws.Range("H78").Value = '{''label'': ''LABEL_0'', ''score'': 0.6986395120620728}'
$ws.Range("H79").Value = '{''label'': ''LABEL_0'', ''score'': 0.7201886773109436}'
$ws.Range("H80").Value = '{''label'': ''LABEL_0'', ''score'': 0.7090784311294556}'
$ws.Range("H81").Value = '{''label'': ''LABEL_0'', ''score'': 0.7135809659957886}'
$ws.Range("H82").Value = '{''label'': ''LABEL_0'', ''score'': 0.6985413432121277}'
$ws.Range("H83").Value = '{''label'': ''LABEL_0'', ''score'': 0.7058203816413879}'
$ws.Range("H84").Value = '{''label'': ''LABEL_0'', ''score'': 0.7120464444160461}'
$ws.Range("H85").Value = '{''label'': ''LABEL_0'', ''score'': 0.71199631690979}'
$ws.Range("H86").Value = '{''label'': ''LABEL_0'', ''score'': 0.7189050316810608}'
$ws.Range("H87").Value = '{''label'': ''LABEL_0'', ''score'': 0.6927416324615479}'
$ws.Range("H88").Value = '{''label'': ''LABEL_0'', ''score'': 0.713079571723938}'
$ws.Range("H89").Value = '{''label'': ''LABEL_0'', ''score'': 0.6928613781929016}'
$ws.Range("H90").Value = '{''label'': ''LABEL_0'', ''score'': 0.7136166095733643}'
$ws.Range("H91").Value = '{''label'': ''LABEL_0'', ''score'': 0.7095724940299988}'
$ws.Range("H92").Value = '{''label'': ''LABEL_0'', ''score'': 0.7074331641197205}'
$ws.Range("H93").Value = '{''label'': ''LABEL_0'', ''score'': 0.6991096138954163}'
$ws.Range("H94").Value = '{''label'': ''LABEL_0'', ''score'': 0.7003313899040222}'
$ws.Range("H95").Value = '{''label'': ''LABEL_0'', ''score'': 0.7140833139419556}'
$ws.Range("H96").Value = '{''label'': ''LABEL_0'', ''score'': 0.7046642899513245}'
$ws.Range("H97").Value = '{''label'': ''LABEL_0'', ''score'': 0.6973547339439392}'
$ws.Range("H98").Value = '{''label'': ''LABEL_0'', ''score'': 0.7223910093307495}'
$ws.Range("H99").Value = '{''label'': ''LABEL_0'', ''score'': 0.6902124881744385}'
$ws.Range("H100").Value = '{''label'': ''LABEL_0'', ''score'': 0.7040569186210632}'
$ws.Range("H101").Value = '{''label'': ''LABEL_0'', ''score'': 0.6902124881744385}'
